$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, $CellRef, $TextValue)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $TextValue
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '67.631.69'
$ws.Range('E2').Value = '  +7.64%  '
Set-TextCell $ws 'D3' '3.543.06'
$ws.Range('E3').Value = '  +10.09%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextCell $ws 'D5' '192.10'
$ws.Range('E5').Value = '  +10.43%  '
Set-TextCell $ws 'D6' '560.68'
$ws.Range('E6').Value = '  +9.19%  '
Set-TextCell $ws 'D7' '3.534.91'
$ws.Range('E7').Value = '  +10.00%  '
$ws.Range('E9').Value = '  -0.08%  '
Set-TextCell $ws 'D10' '0.642'
$ws.Range('E10').Value = '  +7.40%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D11' '56.33'
$ws.Range('E11').Value = '  +7.39%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws 'D12' '0.151'
$ws.Range('E12').Value = '  +16.02%  '
Set-TextCell $ws 'D13' '0.0000274'
$ws.Range('E13').Value = '  +9.17%  '
Set-TextCell $ws 'D14' '9.50'
$ws.Range('E14').Value = '  +7.12%  '
Set-TextCell $ws 'D15' '4.112.96'
$ws.Range('E15').Value = '  +9.91%  '
Set-TextCell $ws 'D16' '3.551.94'
$ws.Range('E16').Value = '  +10.22%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D17' '67.726.40'
$ws.Range('E17').Value = '  +7.84%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D18' '0.122'
$ws.Range('E18').Value = '  +5.51%  '
Set-TextCell $ws 'D19' '18.41'
$ws.Range('E19').Value = '  +7.68%  '
Set-TextCell $ws 'D20' '11.93'
$ws.Range('E20').Value = '  +9.23%  '
$ws.Range('E21').Value = '  +5.00%  '
Set-TextCell $ws 'D22' '407.59'
$ws.Range('E22').Value = '  +11.74%  '
$ws.Range('E23').Value = '  +8.33%  '
Set-TextCell $ws 'D24' '85.52'
$ws.Range('E24').Value = '  +6.99%  '
Set-TextCell $ws 'D25' '4.22'
$ws.Range('E25').Value = '  +8.41%  '
Set-TextCell $ws 'D26' '11.40'
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('E29').Value = '  +7.55%  '
Set-TextCell $ws 'D30' '8.84'
$ws.Range('E30').Value = '  +8.68%  '
Set-TextCell $ws 'D31' '30.61'
$ws.Range('E31').Value = '  +8.97%  '
Set-TextCell $ws 'D32' '689.14'
$ws.Range('E32').Value = '  +5.91%  '
Set-TextCell $ws 'D33' '6.83'
$ws.Range('E33').Value = '  +8.67%  '
Set-TextCell $ws 'D34' '11.84'
$ws.Range('E34').Value = '  +6.97%  '
$ws.Range('E35').Value = '  +8.50%  '
Set-TextCell $ws 'D36' '60.76'
$ws.Range('E36').Value = '  +6.16%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D37' '39.15'
$ws.Range('E37').Value = '  +7.30%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D38' '0.0₃0830'
$ws.Range('E38').Value = '  +18.30%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('E40').Value = '  +7.40%  '
$ws.Range('E41').Value = '  +14.79%  '
$ws.Range('E42').Value = '  +19.63%  '
Set-TextCell $ws 'D43' '1.00'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell $ws 'D44' '3.01'
$ws.Range('E44').Value = '  +16.09%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D45' '3.045.15'
$ws.Range('E45').Value = '  +6.60%  '
$ws.Range('E46').Value = '  +7.95%  '
Set-TextCell $ws 'D47' '0.0423'
$ws.Range('E47').Value = '  +8.72%  '
Set-TextCell $ws 'D48' '3.26'
$ws.Range('E48').Value = '  +11.81%  '
$ws.Range('E49').Value = '  +19.66%  '
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('E51').Value = '  +7.01%  '

Write-Host "Applied 97 cell changes"
